$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 535, shifting existing rows 535-587 down to 537-589
$ws.Rows.Item(535).Insert()
$ws.Rows.Item(536).Insert()

# Populate new row 535
$ws.Cells.Item(535, 1).Value = 10
$ws.Cells.Item(535, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(535, 3).Value = "La Araucanía"
$ws.Cells.Item(535, 4).Value = 44946
$ws.Cells.Item(535, 5).Value = 9
$ws.Cells.Item(535, 6).Value = 100112023
$ws.Cells.Item(535, 7).Value = "Brócoli"
$ws.Cells.Item(535, 8).Value = "Sin especificar"
$ws.Cells.Item(535, 9).Value = "Primera"
$ws.Cells.Item(535, 10).Value = 1450
$ws.Cells.Item(535, 11).Value = 800
$ws.Cells.Item(535, 12).Value = 1000
$ws.Cells.Item(535, 13).Value = 903
$ws.Cells.Item(535, 14).Value = "`$/unidad"
$ws.Cells.Item(535, 15).Value = "Región Metropolitana"
$ws.Cells.Item(535, 16).Value = 903
$ws.Cells.Item(535, 17).Value = 1
$ws.Cells.Item(535, 18).Value = "Hortaliza"

# Populate new row 536
$ws.Cells.Item(536, 1).Value = 10
$ws.Cells.Item(536, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(536, 3).Value = "La Araucanía"
$ws.Cells.Item(536, 4).Value = 44946
$ws.Cells.Item(536, 5).Value = 9
$ws.Cells.Item(536, 6).Value = 100112023
$ws.Cells.Item(536, 7).Value = "Brócoli"
$ws.Cells.Item(536, 8).Value = "Sin especificar"
$ws.Cells.Item(536, 9).Value = "Primera"
$ws.Cells.Item(536, 10).Value = 1250
$ws.Cells.Item(536, 11).Value = 1000
$ws.Cells.Item(536, 12).Value = 1000
$ws.Cells.Item(536, 13).Value = 1000
$ws.Cells.Item(536, 14).Value = "`$/unidad"
$ws.Cells.Item(536, 15).Value = "Región del Maule"
$ws.Cells.Item(536, 16).Value = 1000
$ws.Cells.Item(536, 17).Value = 1
$ws.Cells.Item(536, 18).Value = "Hortaliza"
